$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the two rows (GoalTime, GoalDistance) that will be removed, then delete them.
$ws.Rows("19:20").Select() | Out-Null
$ws.Rows("19:20").Delete() | Out-Null
